$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 29.20950566666667
$ws.Cells.Item(2, 8).Value = 87.628517
$ws.Cells.Item(2, 9).Value = 0.01829497698069002
$ws.Cells.Item(2, 10).Value = 0.01840828041918582
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.506715
$ws.Cells.Item(2, 14).Value = 1.520145
$ws.Cells.Item(2, 15).Value = 0.003122343715987576
$ws.Cells.Item(2, 16).Value = 0.003132472094339857
$ws.Cells.Item(2, 17).Value = 14.800894663885
$ws.Cells.Item(2, 18).Value = 133.208051974965
$ws.Cells.Item(2, 19).Value = 0.00005712320640979484
$ws.Cells.Item(2, 20).Value = 0.00005766342471788238
$ws.Cells.Item(3, 7).Value = 29.20950566666667
$ws.Cells.Item(3, 8).Value = 87.628517
$ws.Cells.Item(3, 9).Value = 0.01829497698069002
$ws.Cells.Item(3, 10).Value = 0.01840828041918582
$ws.Cells.Item(3, 13).Value = 88.13219433333332
$ws.Cells.Item(3, 14).Value = 264.396583
$ws.Cells.Item(3, 15).Value = 0.5430646480820168
$ws.Cells.Item(3, 16).Value = 0.5448262620252092
$ws.Cells.Item(3, 17).Value = 2574.297829795268
$ws.Cells.Item(3, 18).Value = 23168.68046815741
$ws.Cells.Item(3, 19).Value = 0.009935355235687024
$ws.Cells.Item(3, 20).Value = 0.01002931461109686
$ws.Cells.Item(4, 7).Value = 29.20950566666667
$ws.Cells.Item(4, 8).Value = 87.628517
$ws.Cells.Item(4, 9).Value = 0.01829497698069002
$ws.Cells.Item(4, 10).Value = 0.01840828041918582
$ws.Cells.Item(4, 13).Value = 1.5741895
$ws.Cells.Item(4, 14).Value = 3.148379
$ws.Cells.Item(4, 15).Value = 0.009700049718478087
$ws.Cells.Item(4, 16).Value = 0.006487676741301404
$ws.Cells.Item(4, 17).Value = 45.98129712065717
$ws.Cells.Item(4, 18).Value = 275.887782723943
$ws.Cells.Item(4, 19).Value = 0.0001774621863111053
$ws.Cells.Item(4, 20).Value = 0.0001194269727229059
$ws.Cells.Item(5, 7).Value = 29.20950566666667
$ws.Cells.Item(5, 8).Value = 87.628517
$ws.Cells.Item(5, 9).Value = 0.01829497698069002
$ws.Cells.Item(5, 10).Value = 0.01840828041918582
$ws.Cells.Item(5, 13).Value = 72.07364666666666
$ws.Cells.Item(5, 14).Value = 216.22094
$ws.Cells.Item(5, 15).Value = 0.4441129584835175
$ws.Cells.Item(5, 16).Value = 0.4455535891391496
$ws.Cells.Item(5, 17).Value = 2105.235590727331
$ws.Cells.Item(5, 18).Value = 18947.12031654598
$ws.Cells.Item(5, 19).Value = 0.008125036352282095
$ws.Cells.Item(5, 20).Value = 0.008201875410648171
$ws.Cells.Item(6, 9).Value = 0.913374480506715
$ws.Cells.Item(6, 10).Value = 0.9190311407684336
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.506715
$ws.Cells.Item(6, 14).Value = 1.520145
$ws.Cells.Item(6, 15).Value = 0.003122343715987576
$ws.Cells.Item(6, 16).Value = 0.003132472094339857
$ws.Cells.Item(6, 17).Value = 738.932849652085
$ws.Cells.Item(6, 18).Value = 6650.395646868765
$ws.Cells.Item(6, 19).Value = 0.002851869069553559
$ws.Cells.Item(6, 20).Value = 0.002878839402286443
$ws.Cells.Item(7, 9).Value = 0.913374480506715
$ws.Cells.Item(7, 10).Value = 0.9190311407684336
$ws.Cells.Item(7, 13).Value = 88.13219433333332
$ws.Cells.Item(7, 14).Value = 264.396583
$ws.Cells.Item(7, 15).Value = 0.5430646480820168
$ws.Cells.Item(7, 16).Value = 0.5448262620252092
$ws.Cells.Item(7, 17).Value = 128521.5032213795
$ws.Cells.Item(7, 18).Value = 1156693.528992416
$ws.Cells.Item(7, 19).Value = 0.4960213908234741
$ws.Cells.Item(7, 20).Value = 0.5007123011096296
$ws.Cells.Item(8, 9).Value = 0.913374480506715
$ws.Cells.Item(8, 10).Value = 0.9190311407684336
$ws.Cells.Item(8, 13).Value = 1.5741895
$ws.Cells.Item(8, 14).Value = 3.148379
$ws.Cells.Item(8, 15).Value = 0.009700049718478087
$ws.Cells.Item(8, 16).Value = 0.006487676741301404
$ws.Cells.Item(8, 17).Value = 2295.610615686117
$ws.Cells.Item(8, 18).Value = 13773.6636941167
$ws.Cells.Item(8, 19).Value = 0.00885977787250423
$ws.Cells.Item(8, 20).Value = 0.005962376956495064
$ws.Cells.Item(9, 9).Value = 0.913374480506715
$ws.Cells.Item(9, 10).Value = 0.9190311407684336
$ws.Cells.Item(9, 13).Value = 72.07364666666666
$ws.Cells.Item(9, 14).Value = 216.22094
$ws.Cells.Item(9, 15).Value = 0.4441129584835175
$ws.Cells.Item(9, 16).Value = 0.4455535891391496
$ws.Cells.Item(9, 17).Value = 105103.6285016577
$ws.Cells.Item(9, 18).Value = 945932.6565149195
$ws.Cells.Item(9, 19).Value = 0.405641442741183
$ws.Cells.Item(9, 20).Value = 0.4094776233000226
$ws.Cells.Item(10, 7).Value = 57.98602933333333
$ws.Cells.Item(10, 8).Value = 173.958088
$ws.Cells.Item(10, 9).Value = 0.03631876156896331
$ws.Cells.Item(10, 10).Value = 0.03654368891224535
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.506715
$ws.Cells.Item(10, 14).Value = 1.520145
$ws.Cells.Item(10, 15).Value = 0.003122343715987576
$ws.Cells.Item(10, 16).Value = 0.003132472094339857
$ws.Cells.Item(10, 17).Value = 29.38239085364
$ws.Cells.Item(10, 18).Value = 264.44151768276
$ws.Cells.Item(10, 19).Value = 0.0001133996569573037
$ws.Cells.Item(10, 20).Value = 0.0001144720857418454
$ws.Cells.Item(11, 7).Value = 57.98602933333333
$ws.Cells.Item(11, 8).Value = 173.958088
$ws.Cells.Item(11, 9).Value = 0.03631876156896331
$ws.Cells.Item(11, 10).Value = 0.03654368891224535
$ws.Cells.Item(11, 13).Value = 88.13219433333332
$ws.Cells.Item(11, 14).Value = 264.396583
$ws.Cells.Item(11, 15).Value = 0.5430646480820168
$ws.Cells.Item(11, 16).Value = 0.5448262620252092
$ws.Cells.Item(11, 17).Value = 5110.436005823699
$ws.Cells.Item(11, 18).Value = 45993.92405241329
$ws.Cells.Item(11, 19).Value = 0.01972343547022374
$ws.Cells.Item(11, 20).Value = 0.01990996143067071
$ws.Cells.Item(12, 7).Value = 57.98602933333333
$ws.Cells.Item(12, 8).Value = 173.958088
$ws.Cells.Item(12, 9).Value = 0.03631876156896331
$ws.Cells.Item(12, 10).Value = 0.03654368891224535
$ws.Cells.Item(12, 13).Value = 1.5741895
$ws.Cells.Item(12, 14).Value = 3.148379
$ws.Cells.Item(12, 15).Value = 0.009700049718478087
$ws.Cells.Item(12, 16).Value = 0.006487676741301404
$ws.Cells.Item(12, 17).Value = 91.28099852322532
$ws.Cells.Item(12, 18).Value = 547.6859911393519
$ws.Cells.Item(12, 19).Value = 0.0003522937929324953
$ws.Cells.Item(12, 20).Value = 0.0002370836405973282
$ws.Cells.Item(13, 7).Value = 57.98602933333333
$ws.Cells.Item(13, 8).Value = 173.958088
$ws.Cells.Item(13, 9).Value = 0.03631876156896331
$ws.Cells.Item(13, 10).Value = 0.03654368891224535
$ws.Cells.Item(13, 13).Value = 72.07364666666666
$ws.Cells.Item(13, 14).Value = 216.22094
$ws.Cells.Item(13, 15).Value = 0.4441129584835175
$ws.Cells.Item(13, 16).Value = 0.4455535891391496
$ws.Cells.Item(13, 17).Value = 4179.264589773635
$ws.Cells.Item(13, 18).Value = 37613.38130796271
$ws.Cells.Item(13, 19).Value = 0.01612963264884977
$ws.Cells.Item(13, 20).Value = 0.01628217175523546
$ws.Cells.Item(14, 7).Value = 29.481085
$ws.Cells.Item(14, 8).Value = 58.96217
$ws.Cells.Item(14, 9).Value = 0.01846507700595112
$ws.Cells.Item(14, 10).Value = 0.01238628926567028
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.506715
$ws.Cells.Item(14, 14).Value = 1.520145
$ws.Cells.Item(14, 15).Value = 0.003122343715987576
$ws.Cells.Item(14, 16).Value = 0.003132472094339857
$ws.Cells.Item(14, 17).Value = 14.938507985775
$ws.Cells.Item(14, 18).Value = 89.63104791465001
$ws.Cells.Item(14, 19).Value = 0.00005765431715475819
$ws.Cells.Item(14, 20).Value = 0.00003879970547713346
$ws.Cells.Item(15, 7).Value = 29.481085
$ws.Cells.Item(15, 8).Value = 58.96217
$ws.Cells.Item(15, 9).Value = 0.01846507700595112
$ws.Cells.Item(15, 10).Value = 0.01238628926567028
$ws.Cells.Item(15, 13).Value = 88.13219433333332
$ws.Cells.Item(15, 14).Value = 264.396583
$ws.Cells.Item(15, 15).Value = 0.5430646480820168
$ws.Cells.Item(15, 16).Value = 0.5448262620252092
$ws.Cells.Item(15, 17).Value = 2598.232712377518
$ws.Cells.Item(15, 18).Value = 15589.39627426511
$ws.Cells.Item(15, 19).Value = 0.01002773054604419
$ws.Cells.Item(15, 20).Value = 0.006748375680978111
$ws.Cells.Item(16, 7).Value = 29.481085
$ws.Cells.Item(16, 8).Value = 58.96217
$ws.Cells.Item(16, 9).Value = 0.01846507700595112
$ws.Cells.Item(16, 10).Value = 0.01238628926567028
$ws.Cells.Item(16, 13).Value = 1.5741895
$ws.Cells.Item(16, 14).Value = 3.148379
$ws.Cells.Item(16, 15).Value = 0.009700049718478087
$ws.Cells.Item(16, 16).Value = 0.006487676741301404
$ws.Cells.Item(16, 17).Value = 46.4088144556075
$ws.Cells.Item(16, 18).Value = 185.63525782243
$ws.Cells.Item(16, 19).Value = 0.0001791121650132524
$ws.Cells.Item(16, 20).Value = 0.00008035824077992031
$ws.Cells.Item(17, 7).Value = 29.481085
$ws.Cells.Item(17, 8).Value = 58.96217
$ws.Cells.Item(17, 9).Value = 0.01846507700595112
$ws.Cells.Item(17, 10).Value = 0.01238628926567028
$ws.Cells.Item(17, 13).Value = 72.07364666666666
$ws.Cells.Item(17, 14).Value = 216.22094
$ws.Cells.Item(17, 15).Value = 0.4441129584835175
$ws.Cells.Item(17, 16).Value = 0.4455535891391496
$ws.Cells.Item(17, 17).Value = 2124.809303639966
$ws.Cells.Item(17, 18).Value = 12748.8558218398
$ws.Cells.Item(17, 19).Value = 0.008200579977738924
$ws.Cells.Item(17, 20).Value = 0.005518755638435113
$ws.Cells.Item(18, 7).Value = 21.628479
$ws.Cells.Item(18, 8).Value = 64.885437
$ws.Cells.Item(18, 9).Value = 0.01354670393768061
$ws.Cells.Item(18, 10).Value = 0.01363060063446486
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.506715
$ws.Cells.Item(18, 14).Value = 1.520145
$ws.Cells.Item(18, 15).Value = 0.003122343715987576
$ws.Cells.Item(18, 16).Value = 0.003132472094339857
$ws.Cells.Item(18, 17).Value = 10.959474736485
$ws.Cells.Item(18, 18).Value = 98.635272628365
$ws.Cells.Item(18, 19).Value = 0.00004229746591216121
$ws.Cells.Item(18, 20).Value = 0.00004269747611655233
$ws.Cells.Item(19, 7).Value = 21.628479
$ws.Cells.Item(19, 8).Value = 64.885437
$ws.Cells.Item(19, 9).Value = 0.01354670393768061
$ws.Cells.Item(19, 10).Value = 0.01363060063446486
$ws.Cells.Item(19, 13).Value = 88.13219433333332
$ws.Cells.Item(19, 14).Value = 264.396583
$ws.Cells.Item(19, 15).Value = 0.5430646480820168
$ws.Cells.Item(19, 16).Value = 0.5448262620252092
$ws.Cells.Item(19, 17).Value = 1906.165314362418
$ws.Cells.Item(19, 18).Value = 17155.48782926177
$ws.Cells.Item(19, 19).Value = 0.007356736006587792
$ws.Cells.Item(19, 20).Value = 0.007426309192833936
$ws.Cells.Item(20, 7).Value = 21.628479
$ws.Cells.Item(20, 8).Value = 64.885437
$ws.Cells.Item(20, 9).Value = 0.01354670393768061
$ws.Cells.Item(20, 10).Value = 0.01363060063446486
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 1.5741895
$ws.Cells.Item(20, 14).Value = 3.148379
$ws.Cells.Item(20, 15).Value = 0.009700049718478087
$ws.Cells.Item(20, 16).Value = 0.006487676741301404
$ws.Cells.Item(20, 17).Value = 34.0473245427705
$ws.Cells.Item(20, 18).Value = 204.283947256623
$ws.Cells.Item(20, 19).Value = 0.0001314037017170048
$ws.Cells.Item(20, 20).Value = 0.00008843093070618586
$ws.Cells.Item(21, 7).Value = 21.628479
$ws.Cells.Item(21, 8).Value = 64.885437
$ws.Cells.Item(21, 9).Value = 0.01354670393768061
$ws.Cells.Item(21, 10).Value = 0.01363060063446486
$ws.Cells.Item(21, 13).Value = 72.07364666666666
$ws.Cells.Item(21, 14).Value = 216.22094
$ws.Cells.Item(21, 15).Value = 0.4441129584835175
$ws.Cells.Item(21, 16).Value = 0.4455535891391496
$ws.Cells.Item(21, 17).Value = 1558.84335338342
$ws.Cells.Item(21, 18).Value = 14029.59018045078
$ws.Cells.Item(21, 19).Value = 0.006016266763463651
$ws.Cells.Item(21, 20).Value = 0.006073163034808189
